$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the data columns that changed (B-F, I-N; G/H stay as-is)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per row (row index corresponds to worksheet row 2..25)
$data = @{
    2 = @(1.02,1.039267397319292,1.048778010473266,1.048348900972594,1.060497880972812,1.042597308965509,1.044360338277985,1.051536907729423,1.051108994982577,1.063224497929034,1.018711862520257)
    3 = @(1.02,1.040038008623571,1.049386347405405,1.049023903717047,1.061228715480298,1.042766518807142,1.044776980446705,1.051957885441533,1.051596379904423,1.063769980673741,1.018851179828088)
    4 = @(1.02,1.040537460522552,1.049780710038204,1.049461767109275,1.061702747231142,1.042875262405833,1.045046688904985,1.052230317099569,1.051912158373729,1.064123404093753,1.018941335890041)
    5 = @(1.02,1.040747623177223,1.049946672341823,1.049646104083767,1.061902299170968,1.042920798694263,1.045160100026729,1.052344853286411,1.05204500767752,1.064272091534671,1.018979239017938)
    6 = @(1.02,1.040782921675623,1.049974548184999,1.049677070211708,1.061935820528361,1.042928433886748,1.045179143718789,1.05236408471995,1.052067319245164,1.064297063104962,1.01898560319436)
    7 = @(1.02,1.040540267970079,1.049782926960338,1.049464229213095,1.061705412599794,1.042875871570352,1.045048204210766,1.052231847516317,1.051913933137287,1.064125390439583,1.018941842348153)
    8 = @(1.02,1.039527659092206,1.04898344884394,1.048576794051492,1.060744633984502,1.042654648300196,1.04450112019038,1.051679171820807,1.051273623591821,1.063408750132024,1.018758943405266)
    9 = @(1.02,1.037749642357052,1.047580330962424,1.047021472128055,1.059060395462871,1.042259143400126,1.043538016255707,1.050705589273615,1.050148510704375,1.062149537855126,1.018436739041297)
    10 = @(1.02,1.036568673408371,1.046648848754015,1.045990399300621,1.057943604562121,1.04199170240934,1.042896654404075,1.050056826772403,1.04940067903067,1.061312590707715,1.018222024688316)
    11 = @(1.02,1.036058364071075,1.046246464475223,1.045545337776615,1.057461481659481,1.041875013390717,1.042619123156066,1.04977599373584,1.049077412170721,1.060950805744179,1.018129078394088)
    12 = @(1.02,1.035868973484289,1.046097146417785,1.045380234603417,1.057282620667436,1.041831537714805,1.042516064577004,1.049671694209858,1.048957420635259,1.060816517512835,1.018094558492834)
    13 = @(1.02,1.035909591088219,1.046129169035514,1.045415640143688,1.057320976922828,1.041840869362929,1.042538169661903,1.049694066137403,1.048983155391442,1.060845318473939,1.018101962917289)
    14 = @(1.02,1.036042705689312,1.046234118820166,1.045531685949386,1.057446692433886,1.041871422369398,1.042610603704233,1.049767372002506,1.049067491906337,1.060939703492168,1.018126224873672)
    15 = @(1.02,1.036124743387334,1.046298801173334,1.045603213823024,1.057524179253253,1.04189022959488,1.042655236578739,1.04981254011521,1.049119465619555,1.060997869836293,1.018141174076335)
    16 = @(1.02,1.036602563426021,1.046675573945551,1.046019966247435,1.057975632324299,1.04199942809351,1.042915077210483,1.05007546666675,1.04942214490007,1.06133661437813,1.018228193829072)
    17 = @(1.02,1.03690257203166,1.046912170359828,1.046281760347851,1.058259208082637,1.04206768898971,1.043078118295617,1.050240417497093,1.049612155822367,1.061549266832735,1.018282786465549)
    18 = @(1.02,1.037077663919475,1.047050264919352,1.046434595332744,1.058424753247759,1.042107418940889,1.043173234930247,1.050336638663314,1.04972303885622,1.06167336297993,1.018314631942846)
    19 = @(1.02,1.037137382956986,1.047097367114864,1.046486730930374,1.058481223584189,1.04212095130436,1.04320567017387,1.050369448927943,1.049760856023711,1.061715686606263,1.018325490838616)
    20 = @(1.02,1.036870373388088,1.046886776285589,1.046253658347781,1.058228768566753,1.042060374076647,1.043060623706795,1.05022271898129,1.049591763967304,1.061526445062024,1.018276928928324)
    21 = @(1.02,1.036003502295453,1.046203209694327,1.045497507468578,1.05740966623204,1.041862428912631,1.042589272875251,1.049745784843378,1.049042654567912,1.060911906824536,1.018119080204097)
    22 = @(1.02,1.035459398685597,1.045774266408736,1.045023315898947,1.056895943414981,1.041737208596648,1.042293083722835,1.049446001309123,1.048697895495111,1.060526071721017,1.018019860685205)
    23 = @(1.02,1.035747749208238,1.0460015768321,1.045274576455647,1.057168155555972,1.04180366242861,1.042450082808333,1.049604913768463,1.048880612010956,1.060730557503764,1.018072456169228)
    24 = @(1.02,1.036884922254434,1.046898250491078,1.046266356015413,1.058242522440784,1.042063679635208,1.043068528704475,1.050230716153823,1.049600978005694,1.061536757054057,1.018279575689557)
    25 = @(1.02,1.038208539871805,1.047942386781917,1.047422545534639,1.059494757634907,1.042362059141707,1.04378688312838,1.050957238654911,1.050438990922091,1.062474636758461,1.018520023185612)
}

foreach ($r in $data.Keys) {
    $rowValues = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowValues[$i]
    }
}
